$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for team record, matching the formatting of the
# existing header row (bold text, thin border, centered horizontally,
# aligned to top vertically).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$newHeaders = $ws.Range("AD1:AF1")
$newHeaders.Font.Bold = $true
$newHeaders.HorizontalAlignment = -4108
$newHeaders.VerticalAlignment = -4160
$newHeaders.Borders.LineStyle = 1

# Fill in the team record (Wins/Losses/Ties) for every player row (2-49).
$lastRow = 49
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 76
    $ws.Cells.Item($r, 31).Value = 85
    $ws.Cells.Item($r, 32).Value = 0
}
